$d = $word.ActiveDocument

# 1. Ativação date change
$d.Content.Find.Execute(
    "Ativação: 01/01/2020",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Ativação: 01/01/2025",
    2)

# 2. Objetivos (PT) - add "e dentro da empresa"
$d.Content.Find.Execute(
    "Levar os alunos a vivenciarem de forma mais aprofundada problemas reais da indústria para, em equipes, apresentarem as possíveis soluções, de forma que, com isso, desenvolvam habilidades transversais fundamentais para sua vida profissional, tais como trabalho em equipe, gerenciamento de projetos, pro atividade, ao mesmo tempo em que consolidam o conhecimento adquirido durante o curso.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Levar os alunos a vivenciarem de forma mais aprofundada problemas reais da indústria para, em equipes e dentro da empresa, apresentarem as possíveis soluções, de forma que, com isso, desenvolvam habilidades transversais fundamentais para sua vida profissional, tais como trabalho em equipe, gerenciamento de projetos, pro atividade, ao mesmo tempo em que consolidam o conhecimento adquirido durante o curso.",
    2)

# 3. Objetivos (EN) - full rewrite
$d.Content.Find.Execute(
    "To lead students to experience in-depth real problems of the industry in order to present the possible solutions in teams, so that they develop transversal skills that are fundamental to their professional life, such as teamwork, project management, pro activity, at the same time in which they consolidate the knowledge acquired during the course.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Guide students to experience real industry problems in more depth so that, in teams and within the company, they present possible solutions, so that they develop fundamental transversal skills for their professional life, such as teamwork, management of projects, pro activity, while consolidating the knowledge acquired during the course.",
    2)

# 4. Programa (PT) - expand first item with team/meeting/visits text
$d.Content.Find.Execute(
    "Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes; Inovação Sistemática",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes ocorrerá por meio do trabalho em equipes e de reuniões e visitas didáticas realiadas na empresa (mínimo 3);Inovação Sistemática",
    2)

# 5. Programa (EN) - expand first item
$d.Content.Find.Execute(
    "Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Training and work in teams and Communication – the development of essential skills for working in teams will occur through work in teams and meetings and educational visits carried out in the company (minimum 3); Systematic Innovation",
    2)

# 6. Avaliação Critério (PT) - rewrite
$d.Content.Find.Execute(
    "Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Serão feitas três apresentações pelas equipes realizadas no ambiente físico da empresa parceira, as notas serão compostas pelas avaliações dos tutores da empresa e da Escola.",
    2)
